# This script applies the weekly NYPD CompStat data refresh to the 061 precinct
# workbook: bumps the bulletin Volume/Number and reporting week dates, widens
# two columns whose "bestFit" values grew, and refreshes the crime-statistics
# table (rows 14-28 and row 33) with the newly collected figures, including a
# few cells that flip between the numeric "0"/blank-style text markers
# ("0" / "***.*") and real numbers as data becomes available / unavailable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a cell as a genuine number, copying the number-style (s=14,
# the "#,##0" format used throughout the table) from a template cell that is
# guaranteed to remain numeric so the style index used for the destination
# cell matches the target workbook exactly.
# ---------------------------------------------------------------------------
function Set-NumStyleAndValue {
    param($addr, $val)
    $ws.Range("I15").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = $val
}

# Helper: write a cell as the text placeholder "0" or "***.*", using the
# General/text style (s=13) exactly like the other placeholder cells already
# on the sheet.
function Set-TextPlaceholder {
    param($addr, $text)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range("C14").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# Header: bulletin volume/number and the reporting week dates
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "47"

$ws.Range("C9").Characters(27, 10).Text = "11/17/2025"
$ws.Range("C9").Characters(48, 10).Text = "11/23/2025"

# ---------------------------------------------------------------------------
# Column widths: columns I and J (9 and 10) now best-fit to the same width
# as column H (8)
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
Set-TextPlaceholder "G14" "0"
Set-TextPlaceholder "H14" "***.*"

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-NumStyleAndValue "C15" 1
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 15
$ws.Range("K15").Value = -11.764705882352
$ws.Range("L15").Value = 7.142857142857
$ws.Range("M15").Value = 275
$ws.Range("N15").Value = -31.818181818181

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 6
Set-TextPlaceholder "D16" "0"
Set-TextPlaceholder "E16" "***.*"
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 60
$ws.Range("I16").Value = 132
$ws.Range("K16").Value = 7.317073170731
$ws.Range("L16").Value = -15.923566878980
$ws.Range("M16").Value = -0.751879699248
$ws.Range("N16").Value = -83.643122676579

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 58.823529411764
$ws.Range("I17").Value = 265
$ws.Range("J17").Value = 244
$ws.Range("K17").Value = 8.606557377049
$ws.Range("L17").Value = 20.454545454545
$ws.Range("M17").Value = 167.676767676768
$ws.Range("N17").Value = -7.017543859649

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -44.444444444444
$ws.Range("I18").Value = 136
$ws.Range("J18").Value = 153
$ws.Range("K18").Value = -11.111111111111
$ws.Range("L18").Value = -30.256410256410
$ws.Range("M18").Value = -40.611353711790
$ws.Range("N18").Value = -92.838335966298

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -35.714285714285
$ws.Range("F19").Value = 40
$ws.Range("H19").Value = -9.090909090909
$ws.Range("I19").Value = 522
$ws.Range("J19").Value = 546
$ws.Range("K19").Value = -4.395604395604
$ws.Range("L19").Value = -9.688581314878
$ws.Range("M19").Value = 45.403899721448
$ws.Range("N19").Value = -60.122230710466

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = -37.5
$ws.Range("I20").Value = 210
$ws.Range("J20").Value = 236
$ws.Range("K20").Value = -11.016949152542
$ws.Range("L20").Value = 13.513513513513
$ws.Range("M20").Value = 25
$ws.Range("N20").Value = -92.978936810431

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -6.896551724137
$ws.Range("F21").Value = 110
$ws.Range("G21").Value = 115
$ws.Range("H21").Value = -4.347826086956
$ws.Range("I21").Value = 1282
$ws.Range("J21").Value = 1321
$ws.Range("K21").Value = -2.952308856926
$ws.Range("L21").Value = -5.247597930524
$ws.Range("M21").Value = 28.456913827655
$ws.Range("N21").Value = -82.488731047671

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
Set-NumStyleAndValue "D22" 1
$ws.Range("E22").Copy() | Out-Null
$ws.Range("K15").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = -100
Set-TextPlaceholder "F22" "0"
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 19
$ws.Range("K22").Value = -21.052631578947
$ws.Range("M22").Value = -11.764705882352

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
Set-NumStyleAndValue "C23" 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = -42.857142857142
$ws.Range("I23").Value = 56
$ws.Range("J23").Value = 71
$ws.Range("K23").Value = -21.126760563380
$ws.Range("L23").Value = 1.818181818181
$ws.Range("M23").Value = 115.384615384615

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = 10.714285714285
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 101
$ws.Range("H24").Value = -2.970297029702
$ws.Range("I24").Value = 1018
$ws.Range("J24").Value = 1075
$ws.Range("K24").Value = -5.302325581395
$ws.Range("L24").Value = -10.229276895943
$ws.Range("M24").Value = 3.245436105476

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -54.545454545454
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = -23.255813953488
$ws.Range("I25").Value = 377
$ws.Range("J25").Value = 453
$ws.Range("K25").Value = -16.777041942604
$ws.Range("L25").Value = -22.903885480572

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -18.181818181818
$ws.Range("F26").Value = 44
$ws.Range("G26").Value = 56
$ws.Range("H26").Value = -21.428571428571
$ws.Range("I26").Value = 473
$ws.Range("J26").Value = 442
$ws.Range("K26").Value = 7.013574660633
$ws.Range("L26").Value = 28.532608695652
$ws.Range("M26").Value = 1.502145922746

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
Set-NumStyleAndValue "C27" 1
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 16
$ws.Range("K27").Value = -30.434782608695
$ws.Range("L27").Value = -15.789473684210

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 7
$ws.Range("H28").Value = 133.333333333333
$ws.Range("I28").Value = 49
$ws.Range("K28").Value = 48.484848484848
$ws.Range("L28").Value = 25.641025641025

# ---------------------------------------------------------------------------
# Row 33 - Hate Crimes
# ---------------------------------------------------------------------------
$ws.Range("L33").Value = 16.666666666666
